# Add the "Gross Floor Area" intensity table. It mirrors the Building names
# (col B) and Gross Floor Area figures (col H) already present in the first
# table (rows 3-11), re-listed as its own two-column table starting at B28
# with its own header row, ready to be charted ("Added intensity and bar
# graphs").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("B28").Value = "Building"
$ws.Range("C28").Value = "Gross Floor Area"

# Data rows: building name -> gross floor area (sq ft)
$data = @(
    @("10 rock", 487541),
    @("1230 6th", 808600),
    @("30 rock", 2911536),
    @("45 Rock", 1217115),
    @("610 5th", 111101),
    @("50 Rock", 472505),
    @("620 5th", 147585),
    @("1270 6th", 449291),
    @("1 Rock", 482357),
    @("600 5th", 482358),
    @("1 Rock and 600 5th", 964715)
)

$row = 29
foreach ($item in $data) {
    $ws.Range("B$row").Value = $item[0]
    $ws.Range("C$row").Value = $item[1]
    $row++
}

# Restore the view: scrolled so row 10 is the top-left row, with F29 selected
# (matches where the author left the cursor after adding the table).
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F29").Select() | Out-Null
